# feat: support generate enum types
#
# Adds a new "enum<HERO_TYPE|temp>" typed column (O) describing an enum
# generated from the existing `star`(C) values, plus a new data row (7)
# for hero id 106 / star level "01".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column O: enum type declaration + enum member rows ------------
# Written in this exact order so the shared-string table is built in the
# same first-use sequence as the authored workbook (ONE, THREE, type,
# TWO, enum<HERO_TYPE|temp>, temp).
$ws.Range("O4").Value = "ONE"
$ws.Range("O6").Value = "THREE"
$ws.Range("O3").Value = "type"
$ws.Range("O5").Value = "TWO"
$ws.Range("O2").Value = "enum<HERO_TYPE|temp>"
$ws.Range("O7").Value = "temp"

# --- New data row 7 (mirrors rows 4-6's layout: id/star columns are
# text-formatted so leading zeros survive) -------------------------------
$ws.Range("B7:C7").NumberFormat = "@"

$ws.Range("A7").Value = 106
$ws.Range("B7").Value = "001"
$ws.Range("C7").Value = "01"
$ws.Range("D7").Value = 2
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 2
$ws.Range("H7").Value = 2
$ws.Range("I7").Value = 6
$ws.Range("J7").Value = 2
$ws.Range("K7").Value = 6
$ws.Range("L7").Value = 2

# --- Column widths (best-fit sizing for the now-narrower A:N columns and
# the new, wider O column) ----------------------------------------------
$ws.Columns("A").ColumnWidth = 6.142857142857143
$ws.Columns("B:C").ColumnWidth = 3.5714285714285716
$ws.Columns("D").ColumnWidth = 2.7142857142857144
$ws.Columns("E").ColumnWidth = 3.142857142857143
$ws.Columns("F").ColumnWidth = 2.7142857142857144
$ws.Columns("G").ColumnWidth = 3.4285714285714284
$ws.Columns("H:M").ColumnWidth = 2.7142857142857144
$ws.Columns("N").ColumnWidth = 1.2857142857142858
$ws.Columns("O").ColumnWidth = 21.428571428571427

# --- Selection cursor moves to the newly authored type cell ------------
$ws.Range("O11").Select() | Out-Null
